# Asset allocation update (commit: "did 20051 asset allocation")
# Adds two new rows (funds) to the lookup table on Sheet1 and touches the
# two columns whose on-screen widths changed as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------------
# Row 5: First SeaFront Fund
$ws.Range("A5").Value = ".FSFUND HK Equity"
$ws.Range("B5").Value = "First SeaFront Fund"
$ws.Range("D5").Value = "HK"
$ws.Range("E5").Value = "Fund, Other Funds"

# Row 6: Diversified Income Fund
$ws.Range("A6").Value = "CLFLDIF HK Equity"
$ws.Range("B6").Value = "Diversified Income Fund"
$ws.Range("D6").Value = "HK"
$ws.Range("E6").Value = "Fund, Other Funds"

# --- Column width adjustments (columns widened to fit the new entries) ---
$ws.Columns.Item(1).ColumnWidth = 15.26
$ws.Columns.Item(2).ColumnWidth = 20.17

# --- Selection reflects where the user ended up after entering the data --
$ws.Range("E8").Select()
